# Tabulated the pseudo code: indent the if/then/else block of the
# Algorithm section so the nested lines visually line up under "Begin".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B24").Value = "   begin if"
$ws.Range("B25").Value = "        userNameIn = userName"
$ws.Range("B26").Value = "        then"
$ws.Range("B27").Value = "        return false"
$ws.Range("B28").Value = "   end if"

$ws.Range("B27").Select()
